$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-10 (Name, King Checkout, King Stayover, Queen Checkout, Queen Stayover)
$data = @(
    @("Cheyenne", 11, 1, 6, 0),
    @("Float Board", 4, 0, 0, 0),
    @("Hung Board", 0, 0, 1, 0),
    @("Johana", 15, 0, 2, 0),
    @("Juan Carlos", 5, 0, 7, 0),
    @("Mariana", 4, 1, 11, 0),
    @("Nestor", 3, 1, 3, 0),
    @("Shae", 9, 1, 7, 0),
    @("Tameka", 6, 2, 9, 0)
)

# Rows 9 and 10 are new; copy the formatting from an existing name cell (A2)
# before setting values, so the new cells get the same style (s="1").
$ws.Cells.Item(2, 1).Copy($ws.Cells.Item(9, 1))
$ws.Cells.Item(2, 1).Copy($ws.Cells.Item(10, 1))

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $ws.Cells.Item($row, 5).Value = $entry[4]
    $row++
}
